# Append a new bold "TAGS:" heading paragraph followed by a set of plain
# (non-bold) tag-line paragraphs, right after the existing "Role 2:"
# paragraph at the very end of the document body (i.e. immediately before
# the closing sectPr).
#
# A plain Range.InsertAfter("...`r...") would make the new paragraphs
# inherit the bold run formatting of the preceding "Role 2:" paragraph, so
# instead we build the new paragraphs as a small WordprocessingML fragment
# (wrapped in the standard "WordOpenXML" pkg:package envelope) and insert it
# with Range.InsertXML. That way each new paragraph/run gets exactly the
# formatting we specify: bold for the "TAGS:" heading, and completely plain
# (no rPr at all) for every tag line.

$d = $word.ActiveDocument

$tagLines = @(
    "1 Akluq, Aklut -- Clothing, Possessions",
    "1 Ungangkengaita Auluksarait -- Taking Care of Catch",
    "1 Aipangyaraq, Aipaqsaraq -- Marriage",
    "1 Ayagayaraq -- Travel",
    "1 Calirpagyaraq -- Hard Work",
    "1 Kalukat -- Celebrations, Gatherings",
    "1 Yurarpalriit -- Dance Festivals",
    "1 Kevgiq, Kevgiryaraq -- Messenger Feast",
    "1 Ingulagyagaq -- Ingulak Festival",
    "1 Kass'at Tekiteqerraallratni -- Early Western Contact",
    "1 Up'nerkilleq -- Spring Camp",
    "1 Aatailnguut -- Illegitimate Children"
)

function Escape-Xml([string]$text) {
    $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
}

# Bold heading paragraph: "TAGS:"
$bodyXml = "<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>TAGS:</w:t></w:r></w:p>"

# One plain paragraph per tag line -- no pPr/rPr at all.
foreach ($line in $tagLines) {
    $bodyXml += "<w:p><w:r><w:t>$(Escape-Xml $line)</w:t></w:r></w:p>"
}

$fragment = "<?xml version=`"1.0`" standalone=`"yes`"?>" +
    "<?mso-application progid=`"Word.Document`"?>" +
    "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" +
    "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" +
    "<pkg:xmlData>" +
    "<w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">" +
    "<w:body>$bodyXml</w:body>" +
    "</w:document>" +
    "</pkg:xmlData></pkg:part></pkg:package>"

# Collapse to the very end of the document body, then insert the fragment.
$r = $d.Content
$r.Collapse(0)
[void]$r.InsertXML($fragment)
